$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (scenario test-run log rows) ---
$ws.Range("G2").Value = "22/04/2020"

$ws.Range("F3").Value = "jordana"
$ws.Range("G3").Value = "22/04/2020"

$ws.Range("F4").Value = "carlos"
$ws.Range("G4").Value = "22/04/2020"

$ws.Range("A5").Value = "CT 04"
$ws.Range("F5").Value = "jordana"
$ws.Range("G5").Value = "22/04/2020"

$ws.Range("A6").Value = "CT 05"
$ws.Range("B6").Value = "No"
$ws.Range("F6").Value = "carlos"
$ws.Range("G6").Value = "22/04/2020"

# --- Box border around the data table (A2:G6) ---
# Left edge of the table (column A, data rows)
$ws.Range("A2:A6").Borders.Item(7).Weight = -4138
# Right edge of the table (column G, data rows)
$ws.Range("G2:G6").Borders.Item(10).Weight = -4138
# Bottom edge of the table (row 6)
$ws.Range("A6:G6").Borders.Item(9).Weight = -4138

# Bottom row grows slightly to accommodate the thicker border
$ws.Rows.Item(6).RowHeight = 15

# --- Selection moves to G2 ---
[void]$ws.Range("G2").Select()

"done"
